$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 313 (shifts rows 313:403 down to 314:404)
$ws.Rows.Item(313).Insert()

# Populate the newly inserted row 313 with the new data record
$ws.Cells.Item(313, 1).Value = 10
$ws.Cells.Item(313, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(313, 3).Value = "La Araucanía"
$ws.Cells.Item(313, 4).Value = 44551
$ws.Cells.Item(313, 5).Value = 9
$ws.Cells.Item(313, 6).Value = "Fruta"
$ws.Cells.Item(313, 7).Value = 100108
$ws.Cells.Item(313, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(313, 9).Value = 100108006
$ws.Cells.Item(313, 10).Value = "Plátano"
$ws.Cells.Item(313, 11).Value = "Barraganete"
$ws.Cells.Item(313, 12).Value = "Primera"
$ws.Cells.Item(313, 13).Value = 110
$ws.Cells.Item(313, 14).Value = 28000
$ws.Cells.Item(313, 15).Value = 28000
$ws.Cells.Item(313, 16).Value = 28000
$ws.Cells.Item(313, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(313, 18).Value = "Ecuador"
$ws.Cells.Item(313, 19).Value = 1400
$ws.Cells.Item(313, 20).Value = 20
